$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format so numeric-looking strings (e.g. "1.000", "0.9998")
# are preserved exactly as literal text instead of being coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.545.88"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.913.01"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "244.26"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.4854"
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("D8").Value = "0.2891"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "0.06800"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "111.03"
$ws.Range("E10").Value = "  +5.39%  "
$ws.Range("D11").Value = "19.29"
$ws.Range("E11").Value = "  +5.09%  "
$ws.Range("D12").Value = "1.918.98"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "0.07568"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "5.385"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "0.6705"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "295.94"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "30.537.01"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "13.03"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D21").Value = "5.533"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "2.160.91"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "6.448"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "9.466"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "165.99"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "20.30"
$ws.Range("E27").Value = "  -3.49%  "
$ws.Range("D28").Value = "2.076"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "1.433"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").Value = "4.150"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").Value = "4.050"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "0.04986"
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").Value = "0.7350"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").Value = "0.9992"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "0.02036"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").Value = "2.716"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").Value = "2.685"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "2.018"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "109.21"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").Value = "0.4446"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "0.8669"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "5.806"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Value = "0.9996"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "69.55"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").Value = "7.191"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "48.44"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "9.181"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "0.2510"
$ws.Range("E51").Value = "  +0.12%  "

# Row 19 and 20 swap places: ShibaInu now row 19, Dai now row 20, with updated data.
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007591"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  -0.04%  "

# Reset style of D2:E51 back to Normal/default so no stray styling is introduced
# (keeps the cells as plain text cells with no explicit style index, matching the original).
$ws.Range("D2:E51").Style = "Normal"
